# Applies the "Added LU and VD" update to the covid19 cases workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: covid19_cases_switzerland - new case counts for VD (col X)
# and LU (col M), plus a correction to BE (col E) on 2020-03-18.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("covid19_cases_switzerland")

$ws1.Range("X11").Value = 331
$ws1.Range("E12").Value = 131
$ws1.Range("X12").Value = 461
$ws1.Range("M13").Value = 50
$ws1.Range("X13").Value = 695

$ws1.Range("C11").Select()

# ---------------------------------------------------------------------
# Sheet 2: Quellen - new source rows for LU and VD
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Quellen")

$ws2.Range("A15").Value = "LU"
$ws2.Range("B15").Value = "https://www.luzernerzeitung.ch/zentralschweiz/luzern/so-will-die-luzerner-regierung-die-massnahmen-des-bundes-umsetzen-lukb-stellt-50-millionen-franken-bereit-ld.1204954"
$ws2.Range("C15").Value = "'@neph_b"
$ws2.Range("D15").Value = "* Official statement in video"

$ws2.Hyperlinks.Add($ws2.Range("B15"), $ws2.Range("B15").Value) | Out-Null
$ws2.Range("B15").Style = "Hyperlink"

# ---------------------------------------------------------------------
# Sheet 3: demographics - stray "s" keystroke landed in M24
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("demographics")
$ws3.Range("M24").Value = "s"

# ---------------------------------------------------------------------
# back to Sheet 2 to finish the VD row
# ---------------------------------------------------------------------
$ws2.Range("A16").Value = "VD"
$ws2.Range("B16").Value = "https://www.24heures.ch/vaud-regions/Les-contaminations-sont-en-hausse-dans-le-canton-de-Vaud/story/23084946?cache=9efAwefu"
$ws2.Range("C16").Value = "'@f_giroud"

$ws2.Hyperlinks.Add($ws2.Range("B16"), $ws2.Range("B16").Value) | Out-Null
$ws2.Range("B16").Style = "Hyperlink"

# Column B now holds long URLs - widen it to fit (as the author did by hand).
$ws2.Columns.Item(2).AutoFit() | Out-Null

# Touching the print setup materializes the (default) page setup for this
# sheet, matching the paperSize/orientation now present in the file.
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

$ws2.Range("B24").Select()
$ws3.Range("M24").Select()

# Leave the workbook on the data sheet, as in the diff (tabSelected stays
# on covid19_cases_switzerland with C11 highlighted).
$ws1.Select()
$ws1.Range("C11").Select()
